# Handling report mode container
# Update the "Translation" sheet:
#  - Remove the SingleUseId6/7/8 rows and the SingleUseId14/16/17 rows
#    (their containers were re-mapped elsewhere / are no longer needed here)
#  - Append new rows for SingleUseId20..SingleUseId27 (the new report-mode
#    container texts) after the existing SingleUseId19 row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Delete rows bottom-up so earlier row numbers stay valid.
$ws.Rows("16:18").Delete()
$ws.Rows("8:10").Delete()

# After the two deletions the sheet now ends with SingleUseId19 on row 14.
# Append the new rows starting at row 15.
$newRows = @(
  @("SingleUseId20", "Large",  "Right", "LTR", "<value>"),
  @("SingleUseId21", "Large",  "Left",  "LTR", "ODO"),
  @("SingleUseId22", "Medium", "Left",  "LTR", "<value>"),
  @("SingleUseId23", "Medium", "Left",  "LTR", "RANGE"),
  @("SingleUseId24", "Medium", "Right", "LTR", "<value>"),
  @("SingleUseId25", "Medium", "Left",  "LTR", "'000"),
  @("SingleUseId26", "Medium", "Right", "LTR", "<value>"),
  @("SingleUseId27", "Medium", "Left",  "LTR", "KM")
)

$r = 15
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $fcell = $ws.Cells.Item($r, 6)
    $fcell.Value = $row[4]
    # Restore the default (General) cell style in case the leading-apostrophe
    # text marker above changed it, so digit-only values like "000" still
    # round-trip as plain text without leaving a quote-prefix style behind.
    $fcell.Style = "Normal"
    $r = $r + 1
}
